# updated legacy GSC export data
# The "Chart" sheet is a rolling daily export: each refresh drops the oldest
# date row and appends a new date row at the bottom, shifting every other
# row's HTTPS-URL count (column C) up by one day. Column B ("Non-HTTPS
# URLs") stays 0 throughout. We also append the corresponding two new
# trailing dates to cover the extra day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 90
$newLastRow = 91

# --- 1. Snapshot current column C (HTTPS URLs) values for rows 2..lastRow,
#        and the current date text in the first data row, before writing
#        anything (writes below would otherwise clobber what we still need
#        to read).
$cVals = New-Object 'object[]' ($lastRow + 2)
for ($r = 2; $r -le $lastRow; $r++) {
    $cVals[$r] = $ws.Cells.Item($r, 3).Value2
}

$firstDateText = $ws.Cells.Item(2, 1).Text
$firstDate = [DateTime]::ParseExact($firstDateText, "yyyy-MM-dd", $null)

# --- 2. Rewrite rows 2..lastRow: date advances by one day (row r keeps the
#        date that used to belong to row r+1 -> equivalently, start date + (r-2+1) days),
#        column C takes what used to be the NEXT row's value (shift up),
#        column B stays 0.
for ($r = 2; $r -le $lastRow; $r++) {
    $newDate = $firstDate.AddDays($r - 2 + 1)
    $dateText = $newDate.ToString("yyyy-MM-dd")

    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dateText
    $cell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = 0.0

    if ($r -lt $lastRow) {
        $newC = $cVals[$r + 1]
    } else {
        $newC = 0.0
    }
    $ws.Cells.Item($r, 3).Value = $newC
}

# --- 3. Append the new trailing row (91) for the newest date, zero counts.
$newRowDate = $firstDate.AddDays($newLastRow - 2 + 1)
$newRowDateText = $newRowDate.ToString("yyyy-MM-dd")

$newCell = $ws.Cells.Item($newLastRow, 1)
$newCell.NumberFormat = "@"
$newCell.Value = $newRowDateText
$newCell.ClearFormats()

$ws.Cells.Item($newLastRow, 2).Value = 0.0
$ws.Cells.Item($newLastRow, 3).Value = 0.0
